# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet (zh-cn / de-de columns) and on the per-language
#    "Status" column of the zh-cn / de-de detail sheets.
#  - The zh-cn sheet's "Latest Handback DateTime" is refreshed.
#  - The de-de sheet's "Latest Handback DateTime" is refreshed.
#  - Both language sheets' "Error Detail" column is cleared now that the
#    handback succeeded (no more "not the latest" error).
#  - A couple of columns are widened/narrowed to fit the new report layout.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"

# --- Status text: Overview sheet (zh-cn / de-de status columns) ---
$overview.Range("E2").Value2 = $statusText
$overview.Range("F2").Value2 = $statusText

# --- Status text: per-language detail sheets ("Status" column) ---
$zhcn.Range("C2").Value2 = $statusText
$dede.Range("C2").Value2 = $statusText

# --- Latest Handback DateTime refresh ---
$zhcn.Range("K2").Value2 = "2016-10-17 14:40:08"
$dede.Range("K2").Value2 = "2016-10-17 14:40:46"

# --- Error Detail cleared now that handback is in sync ---
$zhcn.Range("P2").Value2 = ""
$dede.Range("P2").Value2 = ""

# --- Column width adjustments to fit the regenerated report ---
$overview.Cells.Item(1, 5).ColumnWidth = 29.09   # column E
$overview.Cells.Item(1, 6).ColumnWidth = 29.09   # column F

$zhcn.Cells.Item(1, 3).ColumnWidth = 29.09    # column C
$zhcn.Cells.Item(1, 16).ColumnWidth = 12.75   # column P

$dede.Cells.Item(1, 3).ColumnWidth = 29.09    # column C
$dede.Cells.Item(1, 16).ColumnWidth = 12.75   # column P
